$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $text) {
    $c = $t.Cell($row, 1)
    $r = $c.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# Simple single-value cell updates (rows 1-12)
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "310"
Set-CellText 5 "0.00001"
Set-CellText 6 "0.00055"
Set-CellText 7 "0.00020"
Set-CellText 9 "0.00034"
Set-CellText 10 "0.00040"
Set-CellText 11 "0.00043"
Set-CellText 12 "0.06957"

# Collapse the multi-run, tab-separated cells (rows 44-46) down to a single value
Set-CellText 44 "99.69"
Set-CellText 45 "0.07"
Set-CellText 46 "22"
